$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $ref, $text)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-CellText $ws "D2" "27.794.60"
Set-CellText $ws "E2" "  +1.06%  "
Set-CellText $ws "D3" "1.900.16"
Set-CellText $ws "E3" "  +2.18%  "
Set-CellText $ws "D4" "1.005"
Set-CellText $ws "E4" "  -0.57%  "
Set-CellText $ws "D5" "315.61"
Set-CellText $ws "E5" "  +1.28%  "
Set-CellText $ws "D6" "1.006"
Set-CellText $ws "E6" "  -0.35%  "
Set-CellText $ws "D7" "0.4822"
Set-CellText $ws "E7" "  +1.15%  "
Set-CellText $ws "D8" "0.3809"
Set-CellText $ws "E8" "  +0.44%  "
Set-CellText $ws "D9" "0.07345"
Set-CellText $ws "E9" "  +0.41%  "
Set-CellText $ws "D10" "0.9257"
Set-CellText $ws "E10" "  -0.33%  "
Set-CellText $ws "D11" "20.63"
Set-CellText $ws "E11" "  -0.07%  "
Set-CellText $ws "B12" "WrappedEther"
Set-CellText $ws "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText $ws "D12" "1.953.57"
Set-CellText $ws "E12" "  +4.77%  "
Set-CellText $ws "B13" "TRON"
Set-CellText $ws "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-CellText $ws "D13" "0.07726"
Set-CellText $ws "E13" "  -0.62%  "
Set-CellText $ws "D14" "5.476"
Set-CellText $ws "E14" "  +0.86%  "
Set-CellText $ws "D15" "6.602"
Set-CellText $ws "E15" "  +0.58%  "
Set-CellText $ws "D16" "91.58"
Set-CellText $ws "E16" "  +1.59%  "
Set-CellText $ws "D17" "1.005"
Set-CellText $ws "E17" "  -0.66%  "
Set-CellText $ws "D18" "0.000008823"
Set-CellText $ws "E18" "  +0.24%  "
Set-CellText $ws "D20" "27.850.16"
Set-CellText $ws "E20" "  +0.91%  "
Set-CellText $ws "D21" "14.57"
Set-CellText $ws "E21" "  -0.55%  "
Set-CellText $ws "D22" "5.142"
Set-CellText $ws "E22" "  +1.04%  "
Set-CellText $ws "D23" "2.161.92"
Set-CellText $ws "E23" "  +1.79%  "
Set-CellText $ws "D24" "10.85"
Set-CellText $ws "E24" "  +1.36%  "
Set-CellText $ws "D25" "1.916"
Set-CellText $ws "E25" "  -0.90%  "
Set-CellText $ws "D26" "154.30"
Set-CellText $ws "E26" "  -1.03%  "
Set-CellText $ws "E27" "  -0.04%  "
Set-CellText $ws "D28" "2.131"
Set-CellText $ws "E28" "  +6.27%  "
Set-CellText $ws "D29" "116.74"
Set-CellText $ws "E29" "  +1.21%  "
Set-CellText $ws "D30" "4.945"
Set-CellText $ws "E30" "  +0.01%  "
Set-CellText $ws "E31" "  +1.12%  "
Set-CellText $ws "D32" "3.189"
Set-CellText $ws "E32" "  -4.08%  "
Set-CellText $ws "D33" "1.242"
Set-CellText $ws "E33" "  +3.61%  "
Set-CellText $ws "D34" "0.7647"
Set-CellText $ws "E34" "  +1.62%  "
Set-CellText $ws "D35" "4.649"
Set-CellText $ws "E35" "  +1.59%  "
Set-CellText $ws "D36" "0.02041"
Set-CellText $ws "E36" "  +0.24%  "
Set-CellText $ws "D37" "2.539"
Set-CellText $ws "E37" "  -6.14%  "
Set-CellText $ws "D38" "1.095"
Set-CellText $ws "E38" "  -2.24%  "
Set-CellText $ws "D39" "0.05274"
Set-CellText $ws "E39" "  -0.93%  "
Set-CellText $ws "D40" "2.988"
Set-CellText $ws "E40" "  +0.25%  "
Set-CellText $ws "D41" "0.5444"
Set-CellText $ws "E41" "  -2.53%  "
Set-CellText $ws "D42" "6.948"
Set-CellText $ws "E42" "  -1.30%  "
Set-CellText $ws "D43" "0.1520"
Set-CellText $ws "E43" "  -0.03%  "
Set-CellText $ws "D44" "8.331"
Set-CellText $ws "E44" "  -1.64%  "
Set-CellText $ws "D45" "10.68"
Set-CellText $ws "E45" "  -0.36%  "
Set-CellText $ws "D46" "109.10"
Set-CellText $ws "E46" "  +4.94%  "
Set-CellText $ws "D47" "0.4794"
Set-CellText $ws "E47" "  -1.37%  "
Set-CellText $ws "D48" "1.006"
Set-CellText $ws "E48" "  -0.29%  "
Set-CellText $ws "E49" "  -1.15%  "
Set-CellText $ws "D50" "67.82"
Set-CellText $ws "E50" "  +0.69%  "
Set-CellText $ws "D51" "0.06073"
Set-CellText $ws "E51" "  -0.47%  "
